# Apply updated dSF (column F) values to the data sheet.
# The F column in rows 2-21 needs to be repulled/recalculated with new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 0
    3  = -2
    4  = -1
    5  = 10
    6  = 2
    7  = -1
    8  = -2
    9  = -1
    10 = 2
    11 = 6
    13 = 3
    14 = 2
    15 = 1
    17 = -3
    18 = 1
    19 = -1
    20 = 2
    21 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("F$row").Value = $newValues[$row]
}
